$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values (rows 2-25) to cycle 1,2,3
for ($r = 2; $r -le 25; $r++) {
    $val = (($r - 2) % 3) + 1
    $ws.Cells.Item($r, 2).Value = $val
}

# Update the selection to B5:B25 with active cell B5
$ws.Range("B5:B25").Select()
